$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-27 Wednesday" "2025-08-28 Thursday"

Replace-Text "87×54=" "37×38="
Replace-Text "68×71=" "11×13="
Replace-Text "79×35=" "33×45="
Replace-Text "12×27=" "49×86="
Replace-Text "82×18=" "86×49="
Replace-Text "55×50=" "61×28="
Replace-Text "86×62=" "42×82="
Replace-Text "35×23=" "74×82="
Replace-Text "25×74=" "93×75="
Replace-Text "36×52=" "96×34="
Replace-Text "30×37=" "93×22="
Replace-Text "21×40=" "89×12="
Replace-Text "30×23=" "40×90="
Replace-Text "90×93=" "85×31="
Replace-Text "98×11=" "14×45="
Replace-Text "53×82=" "18×87="
Replace-Text "18×38=" "83×13="
Replace-Text "83×96=" "91×37="
Replace-Text "90×96=" "98×25="
Replace-Text "18×29=" "81×63="
Replace-Text "79×37=" "28×76="
Replace-Text "15×60=" "22×59="
Replace-Text "20×46=" "82×45="
Replace-Text "93×51=" "76×26="
Replace-Text "92×20=" "46×31="
